# Apply updates to set up for finishing the regex updating.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Shift the "Cases" numbering from column B to column A, and ---
# --- introduce a new "ApplicationGroup" numeric column in B.     ---

# Header row (row 2)
$ws.Range("A2").Value = "Cases"
$ws.Range("B2").Value = "ApplicationGroup"

# Column A: case numbers (what used to live in column B)
$caseNumbers = @(1,2,3,5,6,7,8,9,10,11,13,14,15,19,20,21,22,23)
for ($i = 0; $i -lt $caseNumbers.Length; $i++) {
    $row = 3 + $i
    $ws.Cells.Item($row, 1).Value = $caseNumbers[$i]
}

# Column B: new ApplicationGroup values
$appGroup = @(1,2,3,4,5,5,4,2,6,2,1,2,7,2,7,7,2,2)
for ($i = 0; $i -lt $appGroup.Length; $i++) {
    $row = 3 + $i
    $ws.Cells.Item($row, 2).Value = $appGroup[$i]
}

# Row 14 (case 14) now gets fully filled in, matching the pattern of an
# "AVG(Num1, Num2)" case (2 numbers available, distinct count 2).
$ws.Range("C14").Value = 2
$ws.Range("D14").Value = 2
$ws.Range("E14").Value = "AVG(Num1, Num2)"
$ws.Range("F14").Value = "Num1"
$ws.Range("G14").Value = "Num2"

# --- Remove the old scratch / TODO notes that lived in columns I and J ---
$ws.Range("I6").ClearContents()
$ws.Range("I10").ClearContents()
$ws.Range("I11").ClearContents()
$ws.Range("I14").ClearContents()
$ws.Range("I15").ClearContents()
$ws.Range("I16").ClearContents()
$ws.Range("I17").ClearContents()
$ws.Range("J17").ClearContents()

# --- Update the view: selection now sits on A6 ---
$ws.Range("A6").Select()
